# Update workbook per commit "Mise à jour du fichier via Shiny"
$wb = $excel.ActiveWorkbook

# --- pro sheet (sheet1): B2:B26 values + column B width ---
$wsPro = $wb.Worksheets.Item("pro")
$proVals = New-Object 'object[,]' 25,1
$proVals[0,0] = 3952044.1676466661
$proVals[1,0] = 4294786.304562727
$proVals[2,0] = 4826741.8672162853
$proVals[3,0] = 5370030.689804974
$proVals[4,0] = 4760326.689667305
$proVals[5,0] = 4813983.07266118
$proVals[6,0] = 4172861.8119076509
$proVals[7,0] = 3670255.1596986833
$proVals[8,0] = 3305657.1266281209
$proVals[9,0] = 2783778.3965436206
$proVals[10,0] = 2762758.9136465727
$proVals[11,0] = 2652629.3266975121
$proVals[12,0] = 2515706.0771908681
$proVals[13,0] = 2688116.4137672256
$proVals[14,0] = 2801704.2516777464
$proVals[15,0] = 2478083.5272700186
$proVals[16,0] = 2697868.3770878292
$proVals[17,0] = 3135121.3040956915
$proVals[18,0] = 3417985.2298765611
$proVals[19,0] = 4082449
$proVals[20,0] = 4453107
$proVals[21,0] = 4650557.5002611959
$proVals[22,0] = 4645962.4532787791
$proVals[23,0] = 4904277.9656810788
$proVals[24,0] = 4806192.4063674575
$wsPro.Range("B2:B26").Value = $proVals
$wsPro.Columns.Item(2).ColumnWidth = 14.8

# --- conso sheet (sheet4): B2:B26 values ---
$wsConso = $wb.Worksheets.Item("conso")
$consoVals = New-Object 'object[,]' 25,1
$consoVals[0,0] = 1309378.0682625766
$consoVals[1,0] = 1422934.0956914653
$consoVals[2,0] = 1599182.9152086813
$consoVals[3,0] = 1779180.6252660013
$consoVals[4,0] = 1577172.0753051562
$consoVals[5,0] = 1594951.037141341
$consoVals[6,0] = 1382542.1736807246
$consoVals[7,0] = 1216016.4998023026
$consoVals[8,0] = 1095218.1393083099
$consoVals[9,0] = 922315.12326622615
$consoVals[10,0] = 915350.11792588583
$consoVals[11,0] = 878862.87360652385
$consoVals[12,0] = 833497.71521367901
$consoVals[13,0] = 890618.62564227881
$consoVals[14,0] = 928252.28188350599
$consoVals[15,0] = 821032.28631638025
$consoVals[16,0] = 893852.57149101794
$consoVals[17,0] = 1038722.2998205186
$consoVals[18,0] = 1132440.7242172793
$consoVals[19,0] = 1352587
$consoVals[20,0] = 1428424.9999999998
$consoVals[21,0] = 2046849.1130494985
$consoVals[22,0] = 2030054.2073679301
$consoVals[23,0] = 2141442.2479180484
$consoVals[24,0] = 2098613.4029596876
$wsConso.Range("B2:B26").Value = $consoVals

# --- ind sheet (sheet2): B2:B101 values ---
$wsInd = $wb.Worksheets.Item("ind")
$indVals = New-Object 'object[,]' 100,1
$indVals[0,0] = 1125023.3425613658
$indVals[1,0] = 1090226.8083544241
$indVals[2,0] = 1160122.3252882136
$indVals[3,0] = 1187999.8424038079
$indVals[4,0] = 1195705.4459380703
$indVals[5,0] = 1233814.846493572
$indVals[6,0] = 1284784.2549444556
$indVals[7,0] = 1244827.5131969959
$indVals[8,0] = 1359857.5855413312
$indVals[9,0] = 1337703.1473991002
$indVals[10,0] = 1379778.3788386625
$indVals[11,0] = 1496034.8920433859
$indVals[12,0] = 1467803.8803014997
$indVals[13,0] = 1395009.405350667
$indVals[14,0] = 1837792.0959526207
$indVals[15,0] = 1500096.9284918008
$indVals[16,0] = 1424286.577469293
$indVals[17,0] = 1482369.7913262129
$indVals[18,0] = 1365931.5856525174
$indVals[19,0] = 1224097.3359123014
$indVals[20,0] = 1342546.5641681754
$indVals[21,0] = 1451626.8007279362
$indVals[22,0] = 1439143.7093396643
$indVals[23,0] = 1325324.5208833257
$indVals[24,0] = 1377534.6375781954
$indVals[25,0] = 1336933.9640927704
$indVals[26,0] = 1215159.1550681093
$indVals[27,0] = 888719.72979071073
$indVals[28,0] = 947154.22752379789
$indVals[29,0] = 1053359.6610000415
$indVals[30,0] = 1053601.3131438072
$indVals[31,0] = 1183879.1343040287
$indVals[32,0] = 886161.86016754771
$indVals[33,0] = 1018336.5156854927
$indVals[34,0] = 1042196.1005918888
$indVals[35,0] = 870303.40782873635
$indVals[36,0] = 741955.56728148379
$indVals[37,0] = 837744.3335764159
$indVals[38,0] = 847787.94528842263
$indVals[39,0] = 786903.67809874588
$indVals[40,0] = 802072.76943230117
$indVals[41,0] = 851362.72455269506
$indVals[42,0] = 786784.15724300384
$indVals[43,0] = 749900.9584384599
$indVals[44,0] = 744586.93414252496
$indVals[45,0] = 796669.17300547229
$indVals[46,0] = 752916.74416807014
$indVals[47,0] = 768782.60350286472
$indVals[48,0] = 747866.17004148848
$indVals[49,0] = 769966.60854473023
$indVals[50,0] = 726719.77704931854
$indVals[51,0] = 660299.46205897629
$indVals[52,0] = 672838.28628779622
$indVals[53,0] = 788275.32869607979
$indVals[54,0] = 813109.66628122586
$indVals[55,0] = 829708.63634938269
$indVals[56,0] = 787505.73799767962
$indVals[57,0] = 809796.20612291794
$indVals[58,0] = 816430.04867302976
$indVals[59,0] = 753608.94974844484
$indVals[60,0] = 450602.10558270244
$indVals[61,0] = 506079.1608282841
$indVals[62,0] = 737443.43631750613
$indVals[63,0] = 777621.5707574063
$indVals[64,0] = 829651.77934912872
$indVals[65,0] = 862380.27750158904
$indVals[66,0] = 811096.74753811222
$indVals[67,0] = 859277.53055700765
$indVals[68,0] = 869152.25503955118
$indVals[69,0] = 909758.4638087136
$indVals[70,0] = 888588.54630121484
$indVals[71,0] = 908862.42610672722
$indVals[72,0] = 904840.55841901805
$indVals[73,0] = 911840.31930494553
$indVals[74,0] = 927880.6488524843
$indVals[75,0] = 894427.61105192313
$indVals[76,0] = 988728.01557673444
$indVals[77,0] = 1021501.1402318884
$indVals[78,0] = 997463.90104820009
$indVals[79,0] = 992306.94314317708
$indVals[80,0] = 1082714.8880196426
$indVals[81,0] = 1084952.4636289014
$indVals[82,0] = 1012434.4758852927
$indVals[83,0] = 1060034.597970044
$indVals[84,0] = 1099464.8402971295
$indVals[85,0] = 1092538.1688911328
$indVals[86,0] = 1118404.576937367
$indVals[87,0] = 1161378.2554058982
$indVals[88,0] = 1198303.8728400946
$indVals[89,0] = 1218834.1314740493
$indVals[90,0] = 1206416.4205613739
$indVals[91,0] = 1272963.6485008239
$indVals[92,0] = 1288033.5216520799
$indVals[93,0] = 1304014.3781267575
$indVals[94,0] = 1353567.6410657498
$indVals[95,0] = 1305826.9349690753
$indVals[96,0] = 1203341.8631967304
$indVals[97,0] = 890784.80513578327
$indVals[98,0] = 1048299.1538683418
$indVals[99,0] = 1221523.9630465698
$wsInd.Range("B2:B101").Value = $indVals

# --- VA sheet (sheet3) recomputes automatically via formula =pro!B-conso!B ---

# --- Selections: set activeCell to C105 on every sheet; scroll "ind" to row 75 ---
$wsInd.Select()
$wsInd.Range("C105").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 75

$wsVA = $wb.Worksheets.Item("VA")
$wsVA.Select()
$wsVA.Range("C105").Select()

$wsConso.Select()
$wsConso.Range("C105").Select()

$wsPro.Select()
$wsPro.Range("C105").Select()
